$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.362.95'
$ws.Range("E2").Value = '  -4.39%  '

$ws.Range("D3").Value = '3.678.53'
$ws.Range("E3").Value = '  -5.27%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.94'
$ws.Range("E5").Value = '  -2.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.16'
$ws.Range("E6").Value = '  +4.14%  '

$ws.Range("D7").Value = '3.673.15'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.629'
$ws.Range("E8").Value = '  -6.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.714'
$ws.Range("E10").Value = '  -5.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.160'
$ws.Range("E11").Value = '  -9.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.68'
$ws.Range("E12").Value = '  +2.90%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000292'
$ws.Range("E13").Value = '  -9.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.63'
$ws.Range("E14").Value = '  -7.95%  '

$ws.Range("D15").Value = '4.270.53'
$ws.Range("E15").Value = '  -5.13%  '

$ws.Range("D16").Value = '3.683.01'
$ws.Range("E16").Value = '  -5.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.34'
$ws.Range("E17").Value = '  -8.95%  '

$ws.Range("E18").Value = '  -2.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.82'
$ws.Range("E19").Value = '  -8.34%  '

$ws.Range("E20").Value = '  -7.94%  '

$ws.Range("D21").Value = '68.241.67'
$ws.Range("E21").Value = '  -4.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '410.34'
$ws.Range("E22").Value = '  -7.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.58'
$ws.Range("E23").Value = '  -4.73%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.57'
$ws.Range("E24").Value = '  -6.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.02'
$ws.Range("E25").Value = '  -8.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.67'
$ws.Range("E26").Value = '  -8.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.78'
$ws.Range("E27").Value = '  -9.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.88'
$ws.Range("E28").Value = '  -3.92%  '

$ws.Range("E29").Value = '  +1.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.59'
$ws.Range("E30").Value = '  -8.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.83'
$ws.Range("E31").Value = '  -7.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.35'
$ws.Range("E32").Value = '  -15.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.42'
$ws.Range("E33").Value = '  -8.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.119'
$ws.Range("E34").Value = '  -5.92%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '64.68'
$ws.Range("E35").Value = '  -6.82%  '

$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '602.55'
$ws.Range("E36").Value = '  -5.34%  '

$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '43.24'
$ws.Range("E37").Value = '  -10.65%  '

$ws.Range("D38").Value = '0.0₃0875'
$ws.Range("E38").Value = '  -12.62%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("E40").Value = '  -9.27%  '

$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("E42").Value = '  -6.66%  '

$ws.Range("E43").Value = '  -7.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0441'
$ws.Range("E44").Value = '  -7.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.64'
$ws.Range("E45").Value = '  -8.77%  '

$ws.Range("E46").Value = '  -13.29%  '

$ws.Range("E47").Value = '  -6.68%  '

$ws.Range("E48").Value = '  -6.54%  '

$ws.Range("E49").Value = '  -11.82%  '

$ws.Range("D50").Value = '2.723.88'
$ws.Range("E50").Value = '  -6.50%  '

$ws.Range("E51").Value = '  -5.28%  '
